# The document carried a stray paragraph right after "Jatek" that held
# nothing but a bold space and the (normally hidden) "_GoBack" bookmark.
# That paragraph is removed, and the "_GoBack" bookmark is re-created at
# the end of the "Rendeles" run instead (right after the run, still
# inside its own paragraph).

$d = $word.ActiveDocument

# ---- 1. Delete the stray paragraph that follows "Jatek" -------------------
$jatekRange = $d.Content
$jatekRange.Find.Execute("Játék", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$jatekParaIndex = $jatekRange.Paragraphs.Item(1).Index
$staleParagraph = $d.Paragraphs.Item($jatekParaIndex + 1)
$staleParagraph.Range.Delete()

# ---- 2. Re-create the "_GoBack" bookmark after the "Rendeles" run ---------
# A zero-length Range dropped right on the run/paragraph-mark boundary ends
# up anchored before the run, so instead a unique marker is typed at the
# spot, located again with Find (which yields a well-behaved Range), used
# to plant the bookmark, and then erased.
$marker = "ZzGoBackMarkerZz"

$insertPoint = $d.Content
$insertPoint.Find.Execute("Rendelés", $false, $false, $false, $false, `
                           $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter($marker)

$markerRange = $d.Content
$markerRange.Find.Execute($marker, $false, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Content
$markerRange2.Find.Execute($marker, $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0) | Out-Null
$markerRange2.Text = ""
